$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4868.5435
$ws.Range("I15").Value = 4868.5435
$ws.Range("K15").Value = 14605.6305
$ws.Range("M15").Value = -14436.6305

$ws.Range("H40").Value = 55556492
$ws.Range("I40").Value = 83334136
$ws.Range("K40").Value = 83334136
$ws.Range("M40").Value = -83333961

$ws.Range("H64").Value = 59244.332
$ws.Range("I64").Value = 4439.8
$ws.Range("J64").Value = 127750
$ws.Range("K64").Value = 4439.8
$ws.Range("L64").Value = 127750
$ws.Range("M64").Value = -4191.8
$ws.Range("N64").Value = -128246

$ws.Range("H67").Value = 59244.332
$ws.Range("I67").Value = 4439.8
$ws.Range("J67").Value = 127750
$ws.Range("K67").Value = 4439.8
$ws.Range("L67").Value = 127750
$ws.Range("M67").Value = -3581.8
$ws.Range("N67").Value = -129466

$ws.Range("H106").Value = 2432.55
$ws.Range("I106").Value = 1990.7142
$ws.Range("K106").Value = 1990.7142
$ws.Range("M106").Value = -1359.7142

$ws.Range("H112").Value = 47620870
$ws.Range("J112").Value = 58825628
$ws.Range("L112").Value = 176476884
$ws.Range("N112").Value = -176479100

$ws.Range("H121").Value = 860
$ws.Range("J121").Value = 992.5
$ws.Range("L121").Value = 2977.5
$ws.Range("N121").Value = -6471.5

$ws.Range("H129").Value = 994.125
$ws.Range("I129").Value = 195
$ws.Range("J129").Value = 1108.2858
$ws.Range("K129").Value = 585
$ws.Range("L129").Value = 3324.8574
$ws.Range("M129").Value = 4415
$ws.Range("N129").Value = -13324.8574

$ws.Range("H131").Value = 2113.68
$ws.Range("I131").Value = 1167.5883
$ws.Range("J131").Value = 4124.125
$ws.Range("K131").Value = 3502.7649
$ws.Range("L131").Value = 12372.375
$ws.Range("M131").Value = 1537.2351
$ws.Range("N131").Value = -22452.375

$ws.Range("H132").Value = 148829.31
$ws.Range("I132").Value = 160571.55
$ws.Range("J132").Value = 877.2
$ws.Range("K132").Value = 481714.65
$ws.Range("L132").Value = 2631.6
$ws.Range("M132").Value = -479184.65
$ws.Range("N132").Value = -7691.6

$ws.Range("H137").Value = 45456972
$ws.Range("I137").Value = 1530.6666
$ws.Range("J137").Value = 250006450
$ws.Range("K137").Value = 4591.9998
$ws.Range("L137").Value = 750019350
$ws.Range("M137").Value = -2041.9998
$ws.Range("N137").Value = -750024450

$ws.Range("H138").Value = 1820305.1
$ws.Range("I138").Value = 2440530
$ws.Range("J138").Value = 3932.2144
$ws.Range("K138").Value = 7321590
$ws.Range("L138").Value = 11796.6432
$ws.Range("M138").Value = -7316450
$ws.Range("N138").Value = -22076.6432

$ws.Range("H141").Value = 2082.6924
$ws.Range("I141").Value = 645.2
$ws.Range("J141").Value = 2981.125
$ws.Range("K141").Value = 1935.6
$ws.Range("L141").Value = 8943.375
$ws.Range("M141").Value = 3244.4
$ws.Range("N141").Value = -19303.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5473.8706
$ws.Range("I32").Value = 5567.149
$ws.Range("J32").Value = 4847.5713
$ws.Range("K32").Value = 5567.149
$ws.Range("L32").Value = 4847.5713
$ws.Range("M32").Value = -5280.149
$ws.Range("N32").Value = -5421.5713

$ws.Range("H63").Value = 3790.1
$ws.Range("I63").Value = 2749.375
$ws.Range("J63").Value = 7953
$ws.Range("K63").Value = 2749.375
$ws.Range("L63").Value = 7953
$ws.Range("M63").Value = -2063.375
$ws.Range("N63").Value = -9325

$ws.Range("H66").Value = 3790.1
$ws.Range("I66").Value = 2749.375
$ws.Range("J66").Value = 7953
$ws.Range("K66").Value = 13746.875
$ws.Range("L66").Value = 39765
$ws.Range("M66").Value = -10314.875
$ws.Range("N66").Value = -46629

$ws.Range("H74").Value = 3596.92
$ws.Range("I74").Value = 813.4167
$ws.Range("J74").Value = 10754.5
$ws.Range("K74").Value = 813.4167
$ws.Range("L74").Value = 10754.5
$ws.Range("M74").Value = 60.58330000000001
$ws.Range("N74").Value = -12502.5

$ws.Range("H77").Value = 3596.92
$ws.Range("I77").Value = 813.4167
$ws.Range("J77").Value = 10754.5
$ws.Range("K77").Value = 4067.0835
$ws.Range("L77").Value = 53772.5
$ws.Range("M77").Value = 300.9165000000003
$ws.Range("N77").Value = -62508.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 26315.834
$ws.Range("J35").Value = 26315.834
$ws.Range("L35").Value = 26315.834
$ws.Range("N35").Value = -26935.834

$ws.Range("H82").Value = 16731.75
$ws.Range("J82").Value = 22538.143
$ws.Range("L82").Value = 22538.143
$ws.Range("N82").Value = -23304.143

$ws.Range("H85").Value = 16731.75
$ws.Range("J85").Value = 22538.143
$ws.Range("L85").Value = 22538.143
$ws.Range("N85").Value = -25190.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1291.0714
$ws.Range("I16").Value = 1429.875
$ws.Range("J16").Value = 1106
$ws.Range("K16").Value = 1429.875
$ws.Range("L16").Value = 1106
$ws.Range("M16").Value = -1142.875
$ws.Range("N16").Value = -1680

$ws.Range("H31").Value = 1275.2559
$ws.Range("I31").Value = 933.73334
$ws.Range("J31").Value = 2063.3845
$ws.Range("K31").Value = 933.73334
$ws.Range("L31").Value = 2063.3845
$ws.Range("M31").Value = -638.73334
$ws.Range("N31").Value = -2653.3845

$ws.Range("H34").Value = 1275.2559
$ws.Range("I34").Value = 933.73334
$ws.Range("J34").Value = 2063.3845
$ws.Range("K34").Value = 933.73334
$ws.Range("L34").Value = 2063.3845
$ws.Range("M34").Value = -731.73334
$ws.Range("N34").Value = -2467.3845

$ws.Range("H58").Value = 2585.372
$ws.Range("I58").Value = 996.8095
$ws.Range("J58").Value = 4101.727
$ws.Range("K58").Value = 996.8095
$ws.Range("L58").Value = 4101.727
$ws.Range("M58").Value = -793.8095
$ws.Range("N58").Value = -4507.727

$ws.Range("H105").Value = 715.8570999999999
$ws.Range("I105").Value = 715.8570999999999
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 715.8570999999999
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1031.1429
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 1641.6666
$ws.Range("I107").Value = 2207.5
$ws.Range("K107").Value = 2207.5
$ws.Range("M107").Value = -287.5

$ws.Range("H113").Value = 1291.0714
$ws.Range("I113").Value = 1429.875
$ws.Range("J113").Value = 1106
$ws.Range("K113").Value = 1429.875
$ws.Range("L113").Value = 1106
$ws.Range("M113").Value = 740.125
$ws.Range("N113").Value = -5446

$ws.Range("H132").Value = 1162.2041
$ws.Range("I132").Value = 958.4
$ws.Range("J132").Value = 2068
$ws.Range("K132").Value = 2875.2
$ws.Range("L132").Value = 6204
$ws.Range("M132").Value = -345.1999999999998
$ws.Range("N132").Value = -11264

$ws.Range("H134").Value = 1316.5435
$ws.Range("I134").Value = 1469.6471
$ws.Range("J134").Value = 882.75
$ws.Range("K134").Value = 4408.9413
$ws.Range("L134").Value = 2648.25
$ws.Range("M134").Value = -1873.9413
$ws.Range("N134").Value = -7718.25

$ws.Range("H136").Value = 2585.372
$ws.Range("I136").Value = 996.8095
$ws.Range("J136").Value = 4101.727
$ws.Range("K136").Value = 2990.4285
$ws.Range("L136").Value = 12305.181
$ws.Range("M136").Value = -440.4285
$ws.Range("N136").Value = -17405.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16467047
$ws.Range("I4").Value = 16467047
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 49401141
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -49401029
$ws.Range("N4").ClearContents()

$ws.Range("H86").Value = 2166.6667
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 2166.6667
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 26994.5
$ws.Range("J134").Value = 26994.5
$ws.Range("L134").Value = 80983.5
$ws.Range("N134").Value = -86053.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3732.413
$ws.Range("I132").Value = 4244.448
$ws.Range("J132").Value = 2858.9412
$ws.Range("K132").Value = 12733.344
$ws.Range("L132").Value = 8576.8236
$ws.Range("M132").Value = -10203.344
$ws.Range("N132").Value = -13636.8236

$ws.Range("H136").Value = 2191.2
$ws.Range("I136").Value = 974
$ws.Range("J136").Value = 7060
$ws.Range("K136").Value = 2922
$ws.Range("L136").Value = 21180
$ws.Range("M136").Value = -372
$ws.Range("N136").Value = -26280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4141.6
$ws.Range("I132").Value = 4956.946
$ws.Range("J132").Value = 1821
$ws.Range("K132").Value = 14870.838
$ws.Range("L132").Value = 5463
$ws.Range("M132").Value = -12340.838
$ws.Range("N132").Value = -10523
